$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: P_Subs
$ws.Range("B5").Value = "0.06332391581074477" -as [double]
$ws.Range("C5").Value = "0.06355916260465252" -as [double]
$ws.Range("D5").Value = "0.06382514201487238" -as [double]
$ws.Range("E5").Value = "0.06418211644711853" -as [double]
$ws.Range("F5").Value = "0.053771550193894437" -as [double]

# Row 6: P_ij_1_2
$ws.Range("B6").Value = "0.044494579503654264" -as [double]
$ws.Range("C6").Value = "0.04471003908686566" -as [double]
$ws.Range("D6").Value = "0.04494275373159534" -as [double]
$ws.Range("E6").Value = "0.04529312148249381" -as [double]
$ws.Range("F6").Value = "0.03525584559914957" -as [double]

# Row 7: P_ij_1_3
$ws.Range("B7").Value = "0.009997973752444302" -as [double]
$ws.Range("C7").Value = "0.010357918479412587" -as [double]
$ws.Range("D7").Value = "0.010491231652976257" -as [double]
$ws.Range("E7").Value = "0.010611213653251754" -as [double]
$ws.Range("F7").Value = "0.01042457504513188" -as [double]

# Row 8: P_ij_1_10
$ws.Range("B8").Value = "0.008831362554646203" -as [double]
$ws.Range("C8").Value = "0.008491205038374274" -as [double]
$ws.Range("D8").Value = "0.008391156630300782" -as [double]
$ws.Range("E8").Value = "0.008277781311372977" -as [double]
$ws.Range("F8").Value = "0.008091129549612984" -as [double]

# Row 9: P_ij_2_4
$ws.Range("B9").Value = "0.00500057137543115" -as [double]
$ws.Range("C9").Value = "0.005180603940082182" -as [double]
$ws.Range("D9").Value = "0.005247282880301897" -as [double]
$ws.Range("E9").Value = "0.0053072940254593505" -as [double]
$ws.Range("F9").Value = "0.005213943379754853" -as [double]

# Row 10: P_ij_2_5
$ws.Range("B10").Value = "0.011369152118397572" -as [double]
$ws.Range("C10").Value = "0.011584617583574184" -as [double]
$ws.Range("D10").Value = "0.011728123964554063" -as [double]
$ws.Range("E10").Value = "0.011945571284034685" -as [double]
$ws.Range("F10").Value = "0.005758649387649733" -as [double]

# Row 11: P_ij_2_6
$ws.Range("B11").Value = "0.018121249006797863" -as [double]
$ws.Range("C11").Value = "0.01758111310566857" -as [double]
$ws.Range("D11").Value = "0.017470218260015515" -as [double]
$ws.Range("E11").Value = "0.01742298958043143" -as [double]
$ws.Range("F11").Value = "0.013854953261631183" -as [double]

# Row 12: P_ij_6_7
$ws.Range("B12").Value = "0.0050006681096306" -as [double]
$ws.Range("C12").Value = "0.005180707773371025" -as [double]
$ws.Range("D12").Value = "0.005247389408775227" -as [double]
$ws.Range("E12").Value = "0.005307403011427218" -as [double]
$ws.Range("F12").Value = "0.005214048529859005" -as [double]

# Row 13: P_ij_6_8
$ws.Range("B13").Value = "0.00500063079358504" -as [double]
$ws.Range("C13").Value = "0.005180667720484873" -as [double]
$ws.Range("D13").Value = "0.005247348317387458" -as [double]
$ws.Range("E13").Value = "0.005307360973713176" -as [double]
$ws.Range("F13").Value = "0.005214007963186004" -as [double]

# Row 14: P_ij_8_9
$ws.Range("B14").Value = "0.005000309291350155" -as [double]
$ws.Range("C14").Value = "0.005180322639699659" -as [double]
$ws.Range("D14").Value = "0.005246994289687175" -as [double]
$ws.Range("E14").Value = "0.005306998793222455" -as [double]
$ws.Range("F14").Value = "0.005213658456012855" -as [double]

# Row 15: Q_ij_1_2
$ws.Range("B15").Value = "0.012978671719861955" -as [double]
$ws.Range("C15").Value = "0.014664902722496813" -as [double]
$ws.Range("D15").Value = "0.015422012773512558" -as [double]
$ws.Range("E15").Value = "0.01628853837341492" -as [double]
$ws.Range("F15").Value = "0.015753518592662212" -as [double]

# Row 16: Q_ij_1_3
$ws.Range("B16").Value = "0.005001219967288982" -as [double]
$ws.Range("C16").Value = "0.005181300067762836" -as [double]
$ws.Range("D16").Value = "0.0052479970394988" -as [double]
$ws.Range("E16").Value = "0.005308024609028374" -as [double]
$ws.Range("F16").Value = "0.005214648510507081" -as [double]

# Row 17: Q_ij_1_10
$ws.Range("B17").Value = "0.0026953431676274807" -as [double]
$ws.Range("C17").Value = "0.003271018949616974" -as [double]
$ws.Range("D17").Value = "0.003549848770650065" -as [double]
$ws.Range("E17").Value = "0.003890148426135566" -as [double]
$ws.Range("F17").Value = "0.0036842932696121177" -as [double]

# Row 18: Q_ij_2_4
$ws.Range("B18").Value = "0.0025000759738077485" -as [double]
$ws.Range("C18").Value = "0.0025900908756318476" -as [double]
$ws.Range("D18").Value = "0.0026234299437720363" -as [double]
$ws.Range("E18").Value = "0.0026534352054738727" -as [double]
$ws.Range("F18").Value = "0.002606760378829247" -as [double]

# Row 19: Q_ij_2_5
$ws.Range("B19").Value = "0.001988491321734766" -as [double]
$ws.Range("C19").Value = "0.0021559703211032205" -as [double]
$ws.Range("D19").Value = "0.002218305446441517" -as [double]
$ws.Range("E19").Value = "0.002274100888020736" -as [double]
$ws.Range("F19").Value = "0.0021280208361198543" -as [double]

# Row 20: Q_ij_2_6
$ws.Range("B20").Value = "0.0034773596365156776" -as [double]
$ws.Range("C20").Value = "0.004725704210333591" -as [double]
$ws.Range("D20").Value = "0.005320216020576618" -as [double]
$ws.Range("E20").Value = "0.006040594262941408" -as [double]
$ws.Range("F20").Value = "0.0057964693147216535" -as [double]

# Row 21: Q_ij_6_7
$ws.Range("B21").Value = "0.002500173732131745" -as [double]
$ws.Range("C21").Value = "0.002590195808220757" -as [double]
$ws.Range("D21").Value = "0.002623537600091275" -as [double]
$ws.Range("E21").Value = "0.002653545345321716" -as [double]
$ws.Range("F21").Value = "0.002606866642118722" -as [double]

# Row 22: Q_ij_6_8
$ws.Range("B22").Value = "0.0025001977519201273" -as [double]
$ws.Range("C22").Value = "0.002590221589456038" -as [double]
$ws.Range("D22").Value = "0.002623564049706595" -as [double]
$ws.Range("E22").Value = "0.0026535724039902574" -as [double]
$ws.Range("F22").Value = "0.002606892754028105" -as [double]

# Row 23: Q_ij_8_9
$ws.Range("B23").Value = "0.0024998716494942816" -as [double]
$ws.Range("C23").Value = "0.002589871571107919" -as [double]
$ws.Range("D23").Value = "0.002623204956427134" -as [double]
$ws.Range("E23").Value = "0.0026532050412667596" -as [double]
$ws.Range("F23").Value = "0.0026065382459574075" -as [double]

# Row 24: l_ij_1_2
$ws.Range("B24").Value = "0.002024662168522598" -as [double]
$ws.Range("C24").Value = "0.0020868091920567267" -as [double]
$ws.Range("D24").Value = "0.0021279751529137613" -as [double]
$ws.Range("E24").Value = "0.002183661074367119" -as [double]
$ws.Range("F24").Value = "0.0014055498550281777" -as [double]

# Row 25: l_ij_1_3
$ws.Range("B25").Value = "0.00011779779462322102" -as [double]
$ws.Range("C25").Value = "0.00012643260026237716" -as [double]
$ws.Range("D25").Value = "0.00012970818599585172" -as [double]
$ws.Range("E25").Value = "0.00013269203548402987" -as [double]
$ws.Range("F25").Value = "0.00012806515596165523" -as [double]

# Row 26: l_ij_1_10
$ws.Range("B26").Value = "8.01271931734294e-05" -as [double]
$ws.Range("C26").Value = "7.783592052828577e-05" -as [double]
$ws.Range("D26").Value = "7.804415762771203e-05" -as [double]
$ws.Range("E26").Value = "7.865648042758507e-05" -as [double]
$ws.Range("F26").Value = "7.450178060733178e-05" -as [double]

# Row 27: l_ij_2_4
$ws.Range("B27").Value = "2.9473771198790548e-05" -as [double]
$ws.Range("C27").Value = "3.163491690433745e-05" -as [double]
$ws.Range("D27").Value = "3.245483642782512e-05" -as [double]
$ws.Range("E27").Value = "3.32018397263991e-05" -as [double]
$ws.Range("F27").Value = "3.204208523096305e-05" -as [double]

# Row 28: l_ij_2_5
$ws.Range("B28").Value = "0.00012554357077840337" -as [double]
$ws.Range("C28").Value = "0.00013090660220814636" -as [double]
$ws.Range("D28").Value = "0.00013432468139175897" -as [double]
$ws.Range("E28").Value = "0.00013938119969232358" -as [double]
$ws.Range("F28").Value = "3.554053815717879e-05" -as [double]

# Row 29: l_ij_2_6
$ws.Range("B29").Value = "0.0003209791638092843" -as [double]
$ws.Range("C29").Value = "0.0003124925463587206" -as [double]
$ws.Range("D29").Value = "0.00031447082366101506" -as [double]
$ws.Range("E29").Value = "0.0003206462652343732" -as [double]
$ws.Range("F29").Value = "0.0002126923455470125" -as [double]

# Row 30: l_ij_6_7
$ws.Range("B30").Value = "2.9485639292648346e-05" -as [double]
$ws.Range("C30").Value = "3.164808646761617e-05" -as [double]
$ws.Range("D30").Value = "3.2468630426003706e-05" -as [double]
$ws.Range("E30").Value = "3.3216343344529106e-05" -as [double]
$ws.Range("F30").Value = "3.205404154674865e-05" -as [double]

# Row 31: l_ij_6_8
$ws.Range("B31").Value = "2.948540053924244e-05" -as [double]
$ws.Range("C31").Value = "3.164782096790132e-05" -as [double]
$ws.Range("D31").Value = "3.2468354531195046e-05" -as [double]
$ws.Range("E31").Value = "3.321605786101328e-05" -as [double]
$ws.Range("F31").Value = "3.205377093348442e-05" -as [double]

# Row 32: l_ij_8_9
$ws.Range("B32").Value = "2.9485400539242442e-05" -as [double]
$ws.Range("C32").Value = "3.1647820967901323e-05" -as [double]
$ws.Range("D32").Value = "3.246835453119505e-05" -as [double]
$ws.Range("E32").Value = "3.321605786101327e-05" -as [double]
$ws.Range("F32").Value = "3.205377093348443e-05" -as [double]

# Row 34: v_j_2
$ws.Range("B34").Value = "1.060471486463274" -as [double]
$ws.Range("C34").Value = "1.0604493771525947" -as [double]
$ws.Range("D34").Value = "1.0604386305180635" -as [double]
$ws.Range("E34").Value = "1.0604258243415199" -as [double]
$ws.Range("F34").Value = "1.0604929421723979" -as [double]

# Row 35: v_j_3
$ws.Range("B35").Value = "1.0607372295225344" -as [double]
$ws.Range("C35").Value = "1.060731369416348" -as [double]
$ws.Range("D35").Value = "1.0607291990001404" -as [double]
$ws.Range("E35").Value = "1.0607272456225454" -as [double]
$ws.Range("F35").Value = "1.0607302842086839" -as [double]

# Row 36: v_j_4
$ws.Range("B36").Value = "1.0603071426960828" -as [double]
$ws.Range("C36").Value = "1.0602791167444305" -as [double]
$ws.Range("D36").Value = "1.060266178756688" -as [double]
$ws.Range("E36").Value = "1.0602514003600656" -as [double]
$ws.Range("F36").Value = "1.0603215860883013" -as [double]

# Row 37: v_j_5
$ws.Range("B37").Value = "1.060378363542893" -as [double]
$ws.Range("C37").Value = "1.0603528875405295" -as [double]
$ws.Range("D37").Value = "1.0603405059717077" -as [double]
$ws.Range("E37").Value = "1.0603256995976507" -as [double]
$ws.Range("F37").Value = "1.0604319388490702" -as [double]

# Row 38: v_j_6
$ws.Range("B38").Value = "1.0600940315335743" -as [double]
$ws.Range("C38").Value = "1.0600592674234022" -as [double]
$ws.Range("D38").Value = "1.0600399429619058" -as [double]
$ws.Range("E38").Value = "1.0600152237791725" -as [double]
$ws.Range("F38").Value = "1.0601488649084456" -as [double]

# Row 39: v_j_7
$ws.Range("B39").Value = "1.059880364519082" -as [double]
$ws.Range("C39").Value = "1.0598379079476403" -as [double]
$ws.Range("D39").Value = "1.0598157344184105" -as [double]
$ws.Range("E39").Value = "1.0597884510709572" -as [double]
$ws.Range("F39").Value = "1.0599260809005766" -as [double]

# Row 40: v_j_8
$ws.Range("B40").Value = "1.0599296837797423" -as [double]
$ws.Range("C40").Value = "1.05988900273612" -as [double]
$ws.Range("D40").Value = "1.0598674868104008" -as [double]
$ws.Range("E40").Value = "1.0598407953063955" -as [double]
$ws.Range("F40").Value = "1.0599775044906765" -as [double]

# Row 41: v_j_9
$ws.Range("B41").Value = "1.0598889467340868" -as [double]
$ws.Range("C41").Value = "1.0598467991336291" -as [double]
$ws.Range("D41").Value = "1.0598247400382979" -as [double]
$ws.Range("E41").Value = "1.0597975596814475" -as [double]
$ws.Range("F41").Value = "1.0599350293034706" -as [double]

# Row 42: v_j_10
$ws.Range("B42").Value = "1.0607232105806008" -as [double]
$ws.Range("C42").Value = "1.060719489586395" -as [double]
$ws.Range("D42").Value = "1.0607166981199432" -as [double]
$ws.Range("E42").Value = "1.0607131578303646" -as [double]
$ws.Range("F42").Value = "1.0607192014172477" -as [double]

# Row 43: q_D_j_6
$ws.Range("B43").Value = "0.005036026207556386" -as [double]
$ws.Range("C43").Value = "0.004132669373122418" -as [double]
$ws.Range("D43").Value = "0.003666872680988475" -as [double]
$ws.Range("E43").Value = "0.0030622165506399724" -as [double]
$ws.Range("F43").Value = "0.003089824594449595" -as [double]

# Row 44: q_D_j_10
$ws.Range("B44").Value = "0.0023055277210955417" -as [double]
$ws.Range("C44").Value = "0.0019098431845174946" -as [double]
$ws.Range("D44").Value = "0.0016976849770697743" -as [double]
$ws.Range("E44").Value = "0.0014173930663445257" -as [double]
$ws.Range("F44").Value = "0.0015298780289751922" -as [double]

# Row 45: q_B_j_5
$ws.Range("B45").Value = "0.003012533444693624" -as [double]
$ws.Range("C45").Value = "0.0030250965421749765" -as [double]
$ws.Range("D45").Value = "0.003029452510812138" -as [double]
$ws.Range("E45").Value = "0.00303369127449974" -as [double]
$ws.Range("F45").Value = "0.003085792495030601" -as [double]

# Row 46: q_B_j_6
$ws.Range("B46").Value = "0.0014900733847792735" -as [double]
$ws.Range("C46").Value = "0.0015050655296811194" -as [double]
$ws.Range("D46").Value = "0.0015097221534455907" -as [double]
$ws.Range("E46").Value = "0.0015140737374069323" -as [double]
$ws.Range("F46").Value = "0.00154293987880631" -as [double]

# Row 47: P_c_j_5
$ws.Range("B47").Value = "0.0013713150393736503" -as [double]
$ws.Range("C47").Value = "0.0012268543288419804" -as [double]
$ws.Range("D47").Value = "0.0012370504002260617" -as [double]
$ws.Range("E47").Value = "0.0013345124680544814" -as [double]
$ws.Range("F47").Value = "-8.918345262159885e-09" -as [double]

# Row 48: P_c_j_6
$ws.Range("B48").Value = "0.0004524429015939975" -as [double]
$ws.Range("C48").Value = "0.000592044405570547" -as [double]
$ws.Range("D48").Value = "0.0006810200770422861" -as [double]
$ws.Range("E48").Value = "0.0008602912905415712" -as [double]
$ws.Range("F48").Value = "-8.918550613658096e-09" -as [double]

# Row 49: P_d_j_5
$ws.Range("B49").Value = "4.159894950410076e-08" -as [double]
$ws.Range("C49").Value = "4.1599791410278856e-08" -as [double]
$ws.Range("D49").Value = "4.160017238244907e-08" -as [double]
$ws.Range("E49").Value = "4.16005326142097e-08" -as [double]
$ws.Range("F49").Value = "0.004665508895185284" -as [double]

# Row 50: P_d_j_6
$ws.Range("B50").Value = "4.159337281592296e-08" -as [double]
$ws.Range("C50").Value = "4.1589594244711796e-08" -as [double]
$ws.Range("D50").Value = "4.158872287967354e-08" -as [double]
$ws.Range("E50").Value = "4.158788425774589e-08" -as [double]
$ws.Range("F50").Value = "0.0023335088953943058" -as [double]

# Row 51: B_j_5 (F51 unchanged)
$ws.Range("B51").Value = "0.012966455499037069" -as [double]
$ws.Range("C51").Value = "0.014131923322182834" -as [double]
$ws.Range("D51").Value = "0.015307077412742452" -as [double]
$ws.Range("E51").Value = "0.016574820467359877" -as [double]

# Row 52: B_j_6 (F52 unchanged)
$ws.Range("B52").Value = "0.006263401974016596" -as [double]
$ws.Range("C52").Value = "0.006825800380788358" -as [double]
$ws.Range("D52").Value = "0.0074727256763754985" -as [double]
$ws.Range("E52").Value = "0.00828995862566972" -as [double]

Write-Output "Applied decision_variables update"